# liensMagasinCaissePOD-SCO.xlsx : fill in the missing "-BackOffice-" label
# in B2 (same text as the header-style label already in B1), adjust the
# column B width to fit the new text, and leave the selection on B2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1 already holds " -BackOffice-" ; give B2 the same text.
$ws.Range("B2").Value = $ws.Range("B1").Value2

# Column B grows slightly to accommodate the new text.
$ws.Columns("B:B").ColumnWidth = 19.83

# Selection moves to the cell we just edited.
$ws.Range("B2").Select()
